$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.764.28"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "2.087.79"
$ws.Range("E3").Value = "  +0.55%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.75"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.633"
$ws.Range("E6").Value = "  +1.47%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "58.19"
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.393"
$ws.Range("E9").Value = "  +0.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0780"
$ws.Range("E10").Value = "  -0.29%  "
$ws.Range("E11").Value = "  +3.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.19"
$ws.Range("E12").Value = "  +2.48%  "
$ws.Range("D13").Value = "2.396.49"
$ws.Range("E13").Value = "  +0.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.14"
$ws.Range("E14").Value = "  +1.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.778"
$ws.Range("E15").Value = "  +0.85%  "
$ws.Range("E16").Value = "  +1.17%  "
$ws.Range("D17").Value = "2.087.50"
$ws.Range("E17").Value = "  -0.11%  "
$ws.Range("D18").Value = "37.728.20"
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("E19").Value = "  -1.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.03"
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").Value = "0.0₃0834"
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "229.62"
$ws.Range("E22").Value = "  +0.88%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.39"
$ws.Range("E24").Value = "  -0.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.39"
$ws.Range("E25").Value = "  +0.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.74"
$ws.Range("E26").Value = "  +8.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "171.18"
$ws.Range("E27").Value = "  +1.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.135"
$ws.Range("E28").Value = "  -3.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.53"
$ws.Range("E29").Value = "  +0.47%  "
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.123"
$ws.Range("E31").Value = "  +1.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.69"
$ws.Range("E32").Value = "  +0.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0635"
$ws.Range("E33").Value = "  +1.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.59"
$ws.Range("E34").Value = "  -1.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.49"
$ws.Range("E35").Value = "  +1.06%  "
$ws.Range("E36").Value = "  -0.27%  "
$ws.Range("E37").Value = "  -2.03%  "
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.37"
$ws.Range("E39").Value = "  +0.62%  "
$ws.Range("E40").Value = "  +9.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "101.05"
$ws.Range("E41").Value = "  +3.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0967"
$ws.Range("E42").Value = "  -1.09%  "
$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.91"
$ws.Range("E43").Value = "  +1.43%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.20"
$ws.Range("E44").Value = "  +3.37%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "1.453.23"
$ws.Range("E45").Value = "  +0.25%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.66"
$ws.Range("E46").Value = "  +1.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.14"
$ws.Range("E47").Value = "  -2.31%  "
$ws.Range("E48").Value = "  -0.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.21"
$ws.Range("E49").Value = "  -2.17%  "
$ws.Range("E50").Value = "  -1.86%  "
$ws.Range("D51").Value = "2.279.85"
$ws.Range("E51").Value = "  +0.62%  "
